# Apply the changes described in the commit:
# "Updated references of message, included overdue refID"

$wb = $excel.ActiveWorkbook

$wsNew            = $wb.Worksheets.Item("Transmittals_New")
$wsActionRequired = $wb.Worksheets.Item("Transmittals_New_ActionRequired")
$wsOverdue        = $wb.Worksheets.Item("Transmittals_Overdue")

# ---------------------------------------------------------------------------
# 1. Fix the "message number" formulas on Transmittals_New (sheet1).
#    Each row's M column formula incorrectly referenced N2 for every row;
#    it should reference the N cell of its own row (N3, N4, ... N10).
#    Recalculation will automatically refresh the cached <v> values,
#    including the Delegate-/Reply All- message text for rows 8-10.
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 10; $r++) {
    $mCell = "M$r"
    $nRef  = "N$r"
    $wsNew.Range($mCell).Formula = '=CONCATENATE(ROW()-1," of ",COUNTA(A2:A100)," ",' + $nRef + ')'
}

# ---------------------------------------------------------------------------
# 2. Transmittals_Overdue: update the RefID value and other view tweaks.
# ---------------------------------------------------------------------------
$wsOverdue.Range("A2").Value = "LATFLD-79"

# Column A was resized (no longer auto-fit / bestFit).
$wsOverdue.Columns.Item(1).ColumnWidth = 10.5

# Selection moves back to A2 (top of the data) instead of D11.
$wsOverdue.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Transmittals_New_ActionRequired: scroll the view so column G is visible
#    (selection itself stays on A2).
# ---------------------------------------------------------------------------
$wsActionRequired.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$wsActionRequired.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore Transmittals_New as the active/selected sheet (it was tabSelected
#    before the edit and should remain so).
# ---------------------------------------------------------------------------
$wsNew.Activate() | Out-Null
$wsNew.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Force a full recalculation so all cached formula results (message texts,
#    etc.) are refreshed and consistent with the updated formulas.
# ---------------------------------------------------------------------------
$excel.CalculateFullRebuild()
